# Re-rank the binary-ranking worker table: a handful of rows had their
# identity (prolificid/name/race), re_range position, and realeffort
# score updated as part of recomputing the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 7.311265211180753
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("F3").Value = "Colleen"
$ws.Range("H3").Value = 6.075952185643782
$ws.Range("I3").Value = "White"
$ws.Range("D4").Value = 19
$ws.Range("E4").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("F4").Value = "Jewel"
$ws.Range("H4").Value = 6.068676626552405
$ws.Range("I4").Value = "Black or African American"
$ws.Range("H5").Value = 5.477047804629725
$ws.Range("H6").Value = 5.249471932023906
$ws.Range("H7").Value = 4.260356005502568
$ws.Range("H8").Value = 1.260598627945096
$ws.Range("H9").Value = 1.185192640848691
$ws.Range("D10").Value = 33
$ws.Range("E10").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("F10").Value = "Shaniek"
$ws.Range("H10").Value = 0.3346982378612178
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("F11").Value = "Shadaisia"
$ws.Range("H11").Value = 0.26099946291021
$ws.Range("I11").Value = "Black or African American"
$ws.Range("D12").Value = 32
$ws.Range("E12").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("F12").Value = "Kellie"
$ws.Range("H12").Value = 0.07698541627100014
$ws.Range("I12").Value = "White"
$ws.Range("H13").Value = 0.04919117767745862
$ws.Range("H14").Value = 13.19982871425305
$ws.Range("H15").Value = 8.081433205567341
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("F16").Value = "Matthew"
$ws.Range("H16").Value = 7.048241805590385
$ws.Range("D17").Value = 27
$ws.Range("E17").Value = "5ff8ad350d084e10f500e48a"
$ws.Range("F17").Value = "Drew"
$ws.Range("H17").Value = 7.000947600168775
$ws.Range("D18").Value = 26
$ws.Range("E18").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("F18").Value = "Juan"
$ws.Range("H18").Value = 5.239313832273305
$ws.Range("I18").Value = "Hispanic"
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = "60db4fde6193c50664c9c478"
$ws.Range("F19").Value = "Edosagbe"
$ws.Range("H19").Value = 5.186302527479196
$ws.Range("I19").Value = "Black or African American"
$ws.Range("H20").Value = 5.018452747422359
$ws.Range("H21").Value = 4.155549573790759
$ws.Range("H22").Value = 4.034175108618071
$ws.Range("H23").Value = 3.262000135003892
$ws.Range("H24").Value = 2.30063950806506
$ws.Range("H25").Value = 2.244526961475056
